$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(35,1).Value = 2014
$ws.Cells.Item(35,2).Value = 'Receitas de Exploração de Recursos Naturais'
$ws.Cells.Item(35,3).Value = 'Total (3)'
$ws.Cells.Item(35,4).Value = -1.121216821977611
$ws.Cells.Item(36,1).Value = 2015
$ws.Cells.Item(36,2).Value = 'Receitas de Exploração de Recursos Naturais'
$ws.Cells.Item(36,3).Value = 'Total (3)'
$ws.Cells.Item(36,4).Value = -46.817824810541
$ws.Cells.Item(37,1).Value = 2016
$ws.Cells.Item(37,2).Value = 'Receitas de Exploração de Recursos Naturais'
$ws.Cells.Item(37,3).Value = 'Total (3)'
$ws.Cells.Item(37,4).Value = -31.44547510309983
$ws.Cells.Item(38,1).Value = 2017
$ws.Cells.Item(38,2).Value = 'Receitas de Exploração de Recursos Naturais'
$ws.Cells.Item(38,3).Value = 'Total (3)'
$ws.Cells.Item(38,4).Value = -2.948205694391848
$ws.Cells.Item(39,1).Value = 2018
$ws.Cells.Item(39,2).Value = 'Receitas de Exploração de Recursos Naturais'
$ws.Cells.Item(39,3).Value = 'Total (3)'
$ws.Cells.Item(39,4).Value = 15.01771256682733
$ws.Cells.Item(40,1).Value = 2019
$ws.Cells.Item(40,2).Value = 'Receitas de Exploração de Recursos Naturais'
$ws.Cells.Item(40,3).Value = 'Total (3)'
$ws.Cells.Item(40,4).Value = -20.03643175044765
$ws.Cells.Item(41,1).Value = 2020
$ws.Cells.Item(41,2).Value = 'Receitas de Exploração de Recursos Naturais'
$ws.Cells.Item(41,3).Value = 'Total (3)'
$ws.Cells.Item(41,4).Value = -20.67228986991331
$ws.Cells.Item(42,1).Value = 2021
$ws.Cells.Item(42,2).Value = 'Receitas de Exploração de Recursos Naturais'
$ws.Cells.Item(42,3).Value = 'Total (3)'
$ws.Cells.Item(42,4).Value = 11.91783576722245
$ws.Cells.Item(43,1).Value = 2022
$ws.Cells.Item(43,2).Value = 'Receitas de Exploração de Recursos Naturais'
$ws.Cells.Item(43,3).Value = 'Total (3)'
$ws.Cells.Item(43,4).Value = 44.82592946317314
$ws.Cells.Item(44,1).Value = 2023
$ws.Cells.Item(44,2).Value = 'Receitas de Exploração de Recursos Naturais'
$ws.Cells.Item(44,3).Value = 'Total (3)'
$ws.Cells.Item(44,4).Value = -35.87025361669786
$ws.Cells.Item(45,1).Value = 2024
$ws.Cells.Item(45,2).Value = 'Receitas de Exploração de Recursos Naturais'
$ws.Cells.Item(45,3).Value = 'Total (3)'
$ws.Cells.Item(45,4).Value = 26.44573759185258
$ws.Cells.Item(46,1).Value = 2014
$ws.Cells.Item(46,2).Value = 'Transferências Federais'
$ws.Cells.Item(46,3).Value = 'Fundeb'
$ws.Cells.Item(46,4).Value = -1.050745985120649
$ws.Cells.Item(47,1).Value = 2015
$ws.Cells.Item(47,2).Value = 'Transferências Federais'
$ws.Cells.Item(47,3).Value = 'Fundeb'
$ws.Cells.Item(47,4).Value = -7.493146079046776
$ws.Cells.Item(48,1).Value = 2016
$ws.Cells.Item(48,2).Value = 'Transferências Federais'
$ws.Cells.Item(48,3).Value = 'Fundeb'
$ws.Cells.Item(48,4).Value = 1.677754082062677
$ws.Cells.Item(49,1).Value = 2017
$ws.Cells.Item(49,2).Value = 'Transferências Federais'
$ws.Cells.Item(49,3).Value = 'Fundeb'
$ws.Cells.Item(49,4).Value = -3.335791029940871
$ws.Cells.Item(50,1).Value = 2018
$ws.Cells.Item(50,2).Value = 'Transferências Federais'
$ws.Cells.Item(50,3).Value = 'Fundeb'
$ws.Cells.Item(50,4).Value = 5.808699603470968
$ws.Cells.Item(51,1).Value = 2019
$ws.Cells.Item(51,2).Value = 'Transferências Federais'
$ws.Cells.Item(51,3).Value = 'Fundeb'
$ws.Cells.Item(51,4).Value = -1.74237194029041
$ws.Cells.Item(52,1).Value = 2020
$ws.Cells.Item(52,2).Value = 'Transferências Federais'
$ws.Cells.Item(52,3).Value = 'Fundeb'
$ws.Cells.Item(52,4).Value = -7.79704214867224
$ws.Cells.Item(53,1).Value = 2021
$ws.Cells.Item(53,2).Value = 'Transferências Federais'
$ws.Cells.Item(53,3).Value = 'Fundeb'
$ws.Cells.Item(53,4).Value = 17.82952282634764
$ws.Cells.Item(54,1).Value = 2022
$ws.Cells.Item(54,2).Value = 'Transferências Federais'
$ws.Cells.Item(54,3).Value = 'Fundeb'
$ws.Cells.Item(54,4).Value = 14.13996924257421
$ws.Cells.Item(55,1).Value = 2023
$ws.Cells.Item(55,2).Value = 'Transferências Federais'
$ws.Cells.Item(55,3).Value = 'Fundeb'
$ws.Cells.Item(55,4).Value = 3.254075410574742
$ws.Cells.Item(56,1).Value = 2024
$ws.Cells.Item(56,2).Value = 'Transferências Federais'
$ws.Cells.Item(56,3).Value = 'Fundeb'
$ws.Cells.Item(56,4).Value = 3.251088949624781
$ws.Cells.Item(57,1).Value = 2014
$ws.Cells.Item(57,2).Value = 'Transferências Federais'
$ws.Cells.Item(57,3).Value = 'Fundo de Participação dos Estados'
$ws.Cells.Item(57,4).Value = 2.409290835890476
$ws.Cells.Item(58,1).Value = 2015
$ws.Cells.Item(58,2).Value = 'Transferências Federais'
$ws.Cells.Item(58,3).Value = 'Fundo de Participação dos Estados'
$ws.Cells.Item(58,4).Value = -4.946685733405854
$ws.Cells.Item(59,1).Value = 2016
$ws.Cells.Item(59,2).Value = 'Transferências Federais'
$ws.Cells.Item(59,3).Value = 'Fundo de Participação dos Estados'
$ws.Cells.Item(59,4).Value = 5.961975945722298
$ws.Cells.Item(60,1).Value = 2017
$ws.Cells.Item(60,2).Value = 'Transferências Federais'
$ws.Cells.Item(60,3).Value = 'Fundo de Participação dos Estados'
$ws.Cells.Item(60,4).Value = -6.389518142136009
$ws.Cells.Item(61,1).Value = 2018
$ws.Cells.Item(61,2).Value = 'Transferências Federais'
$ws.Cells.Item(61,3).Value = 'Fundo de Participação dos Estados'
$ws.Cells.Item(61,4).Value = 2.762736106076624
$ws.Cells.Item(62,1).Value = 2019
$ws.Cells.Item(62,2).Value = 'Transferências Federais'
$ws.Cells.Item(62,3).Value = 'Fundo de Participação dos Estados'
$ws.Cells.Item(62,4).Value = 4.889695411025086
$ws.Cells.Item(63,1).Value = 2020
$ws.Cells.Item(63,2).Value = 'Transferências Federais'
$ws.Cells.Item(63,3).Value = 'Fundo de Participação dos Estados'
$ws.Cells.Item(63,4).Value = -8.527616269027838
$ws.Cells.Item(64,1).Value = 2021
$ws.Cells.Item(64,2).Value = 'Transferências Federais'
$ws.Cells.Item(64,3).Value = 'Fundo de Participação dos Estados'
$ws.Cells.Item(64,4).Value = 20.7066685806667
$ws.Cells.Item(65,1).Value = 2022
$ws.Cells.Item(65,2).Value = 'Transferências Federais'
$ws.Cells.Item(65,3).Value = 'Fundo de Participação dos Estados'
$ws.Cells.Item(65,4).Value = 14.95317686156681
$ws.Cells.Item(66,1).Value = 2023
$ws.Cells.Item(66,2).Value = 'Transferências Federais'
$ws.Cells.Item(66,3).Value = 'Fundo de Participação dos Estados'
$ws.Cells.Item(66,4).Value = 0.5246649380495905
$ws.Cells.Item(67,1).Value = 2024
$ws.Cells.Item(67,2).Value = 'Transferências Federais'
$ws.Cells.Item(67,3).Value = 'Fundo de Participação dos Estados'
$ws.Cells.Item(67,4).Value = 9.330859627920308
$ws.Cells.Item(68,1).Value = 2014
$ws.Cells.Item(68,2).Value = 'Transferências Federais'
$ws.Cells.Item(68,3).Value = 'Total (2)'
$ws.Cells.Item(68,4).Value = 0.3913311029190458
$ws.Cells.Item(69,1).Value = 2015
$ws.Cells.Item(69,2).Value = 'Transferências Federais'
$ws.Cells.Item(69,3).Value = 'Total (2)'
$ws.Cells.Item(69,4).Value = -7.525355182452831
$ws.Cells.Item(70,1).Value = 2016
$ws.Cells.Item(70,2).Value = 'Transferências Federais'
$ws.Cells.Item(70,3).Value = 'Total (2)'
$ws.Cells.Item(70,4).Value = 4.714949766020449
$ws.Cells.Item(71,1).Value = 2017
$ws.Cells.Item(71,2).Value = 'Transferências Federais'
$ws.Cells.Item(71,3).Value = 'Total (2)'
$ws.Cells.Item(71,4).Value = -4.934173426721832
$ws.Cells.Item(72,1).Value = 2018
$ws.Cells.Item(72,2).Value = 'Transferências Federais'
$ws.Cells.Item(72,3).Value = 'Total (2)'
$ws.Cells.Item(72,4).Value = 5.774829641842771
$ws.Cells.Item(73,1).Value = 2019
$ws.Cells.Item(73,2).Value = 'Transferências Federais'
$ws.Cells.Item(73,3).Value = 'Total (2)'
$ws.Cells.Item(73,4).Value = 5.802324679288562
$ws.Cells.Item(74,1).Value = 2020
$ws.Cells.Item(74,2).Value = 'Transferências Federais'
$ws.Cells.Item(74,3).Value = 'Total (2)'
$ws.Cells.Item(74,4).Value = 7.581511795768359
$ws.Cells.Item(75,1).Value = 2021
$ws.Cells.Item(75,2).Value = 'Transferências Federais'
$ws.Cells.Item(75,3).Value = 'Total (2)'
$ws.Cells.Item(75,4).Value = -1.876111608518882
$ws.Cells.Item(76,1).Value = 2022
$ws.Cells.Item(76,2).Value = 'Transferências Federais'
$ws.Cells.Item(76,3).Value = 'Total (2)'
$ws.Cells.Item(76,4).Value = 13.31183179031803
$ws.Cells.Item(77,1).Value = 2023
$ws.Cells.Item(77,2).Value = 'Transferências Federais'
$ws.Cells.Item(77,3).Value = 'Total (2)'
$ws.Cells.Item(77,4).Value = 4.645733788933248
$ws.Cells.Item(78,1).Value = 2024
$ws.Cells.Item(78,2).Value = 'Transferências Federais'
$ws.Cells.Item(78,3).Value = 'Total (2)'
$ws.Cells.Item(78,4).Value = 9.290258964394194

$ws.Range("A79:D87").EntireRow.Delete()
